{"js": "// The document stores per-section identifiers as an XML-ish tag:\n//   <id>p081v_a1</id>, <id>p081v_a2</id>, ...\n// each split across three runs: \"<id>\" (Courier New/olive), the bare\n// id text (black, default font), and \"</id>\" (Courier New/olive again).\n// The edit collapses the \"_aN\" placeholder ids into their final form\n// \"_N\" (dropping the \"a\") and merges the three runs for each <id> tag\n// back into a single run (keeping the tag-run's formatting).\n\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n\n// Discover every \"<id>...._aN</id>\" occurrence currently in the body so\n// the same logic works regardless of how many there are / what the\n// numeric prefix is.\nconst idPattern = /<id>([^<>]*?)_a(\\d+)<\\/id>/g;\nconst targets = [];\nlet match;\nwhile ((match = idPattern.exec(body.text)) !== null) {\n  targets.push({\n    oldText: match[0],                       // e.g. \"<id>p081v_a1</id>\"\n    newText: `<id>${match[1]}_${match[2]}</id>`, // e.g. \"<id>p081v_1</id>\"\n  });\n}\n\nfor (const target of targets) {\n  const found = context.document.body.search(target.oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) continue;\n\n  // search() returns the whole \"<id>...</id>\" span as a single range\n  // even though it is backed by several runs in the OOXML; replacing its\n  // text collapses those runs into one run that keeps the formatting of\n  // the first (leading \"<id>\") run, matching the target edit.\n  found.items[0].insertText(target.newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document stores per-section identifiers as an XML-ish tag:\n#   <id>p081v_a1</id>, <id>p081v_a2</id>, ...\n# each split across three runs: \"<id>\" (Courier New/olive), the bare\n# id text (black, default font), and \"</id>\" (Courier New/olive again).\n# The edit collapses the \"_aN\" placeholder ids into their final form\n# \"_N\" (dropping the \"a\") and merges the three runs for each <id> tag\n# back into a single run (keeping the tag-run's formatting).\n\n$d = $word.ActiveDocument\n\n# Discover every \"<id>...._aN</id>\" occurrence currently in the body so\n# the same logic works regardless of how many there are / what the\n# numeric prefix is.\n$bodyText = $d.Content.Text\n$idMatches = [regex]::Matches($bodyText, '<id>([^<>]*?)_a(\\d+)</id>')\n\nforeach ($m in $idMatches) {\n    $oldText = $m.Value                                              # e.g. \"<id>p081v_a1</id>\"\n    $newText = \"<id>\" + $m.Groups[1].Value + \"_\" + $m.Groups[2].Value + \"</id>\"  # e.g. \"<id>p081v_1</id>\"\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWildcards = $false\n\n    # Find.Execute collapses $range to the matched text (the whole\n    # \"<id>...</id>\" span, even though it is backed by several runs in\n    # the OOXML); assigning .Text replaces that span and collapses the\n    # backing runs into one run that keeps the formatting of the first\n    # (leading \"<id>\") run, matching the target edit.\n    $found = $range.Find.Execute($oldText)\n    if ($found) {\n        $range.Text = $newText\n    }\n}\n"}
